$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl13"
$ws.Range("C2").Value = "Ackr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.046465333333333
$ws.Range("H2").Value = 12.139396
$ws.Range("I2").Value = 0.5162107379131895
$ws.Range("J2").Value = 0.5162107379131895
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1143813333333333
$ws.Range("N2").Value = 0.343144
$ws.Range("O2").Value = 0.128300337591142
$ws.Range("P2").Value = 0.1283003375911419
$ws.Range("Q2").Value = 0.4628401001137777
$ws.Range("R2").Value = 4.165560901024
$ws.Range("S2").Value = 0.06623001194243472
$ws.Range("T2").Value = 0.06623001194243471

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cxcl13"
$ws.Range("C3").Value = "Ackr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.046465333333333
$ws.Range("H3").Value = 12.139396
$ws.Range("I3").Value = 0.5162107379131895
$ws.Range("J3").Value = 0.5162107379131895
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7200953333333334
$ws.Range("N3").Value = 2.160286
$ws.Range("O3").Value = 0.8077233554817153
$ws.Range("P3").Value = 0.8077233554817151
$ws.Range("Q3").Value = 2.913840803028445
$ws.Range("R3").Value = 26.224567227256
$ws.Range("S3").Value = 0.4169554693629337
$ws.Range("T3").Value = 0.4169554693629337

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl13"
$ws.Range("C4").Value = "Ackr4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.046465333333333
$ws.Range("H4").Value = 12.139396
$ws.Range("I4").Value = 0.5162107379131895
$ws.Range("J4").Value = 0.5162107379131895
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03357866666666667
$ws.Range("N4").Value = 0.100736
$ws.Range("O4").Value = 0.03766483694187069
$ws.Range("P4").Value = 0.03766483694187069
$ws.Range("Q4").Value = 0.1358749106062222
$ws.Range("R4").Value = 1.222874195456
$ws.Range("S4").Value = 0.01944299327114303
$ws.Range("T4").Value = 0.01944299327114303

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl13"
$ws.Range("C5").Value = "Ackr4"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.046465333333333
$ws.Range("H5").Value = 12.139396
$ws.Range("I5").Value = 0.5162107379131895
$ws.Range("J5").Value = 0.5162107379131895
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.023457
$ws.Range("N5").Value = 0.070371
$ws.Range("O5").Value = 0.02631146998527222
$ws.Range("P5").Value = 0.02631146998527222
$ws.Range("Q5").Value = 0.094917937324
$ws.Range("R5").Value = 0.854261435916
$ws.Range("S5").Value = 0.01358226333667811
$ws.Range("T5").Value = 0.01358226333667811

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Cxcl13"
$ws.Range("C6").Value = "Ackr4"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.155986
$ws.Range("H6").Value = 0.467958
$ws.Range("I6").Value = 0.01989925565426652
$ws.Range("J6").Value = 0.01989925565426652
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1143813333333333
$ws.Range("N6").Value = 0.343144
$ws.Range("O6").Value = 0.128300337591142
$ws.Range("P6").Value = 0.1283003375911419
$ws.Range("Q6").Value = 0.01784188666133333
$ws.Range("R6").Value = 0.160576979952
$ws.Range("S6").Value = 0.002553081218254834
$ws.Range("T6").Value = 0.002553081218254834

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Cxcl13"
$ws.Range("C7").Value = "Ackr4"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.155986
$ws.Range("H7").Value = 0.467958
$ws.Range("I7").Value = 0.01989925565426652
$ws.Range("J7").Value = 0.01989925565426652
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7200953333333334
$ws.Range("N7").Value = 2.160286
$ws.Range("O7").Value = 0.8077233554817153
$ws.Range("P7").Value = 0.8077233554817151
$ws.Range("Q7").Value = 0.1123247906653333
$ws.Range("R7").Value = 1.010923115988
$ws.Range("S7").Value = 0.01607309354865265
$ws.Range("T7").Value = 0.01607309354865264

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Cxcl13"
$ws.Range("C8").Value = "Ackr4"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.155986
$ws.Range("H8").Value = 0.467958
$ws.Range("I8").Value = 0.01989925565426652
$ws.Range("J8").Value = 0.01989925565426652
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03357866666666667
$ws.Range("N8").Value = 0.100736
$ws.Range("O8").Value = 0.03766483694187069
$ws.Range("P8").Value = 0.03766483694187069
$ws.Range("Q8").Value = 0.005237801898666667
$ws.Range("R8").Value = 0.047140217088
$ws.Range("S8").Value = 0.0007495022194825466
$ws.Range("T8").Value = 0.0007495022194825466

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Cxcl13"
$ws.Range("C9").Value = "Ackr4"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.155986
$ws.Range("H9").Value = 0.467958
$ws.Range("I9").Value = 0.01989925565426652
$ws.Range("J9").Value = 0.01989925565426652
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.023457
$ws.Range("N9").Value = 0.070371
$ws.Range("O9").Value = 0.02631146998527222
$ws.Range("P9").Value = 0.02631146998527222
$ws.Range("Q9").Value = 0.003658963602
$ws.Range("R9").Value = 0.032930672418
$ws.Range("S9").Value = 0.000523578667876492
$ws.Range("T9").Value = 0.0005235786678764919

# Row 10
$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Cxcl13"
$ws.Range("C10").Value = "Ackr4"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.636334333333334
$ws.Range("H10").Value = 10.909003
$ws.Range("I10").Value = 0.463890006432544
$ws.Range("J10").Value = 0.463890006432544
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1143813333333333
$ws.Range("N10").Value = 0.343144
$ws.Range("O10").Value = 0.128300337591142
$ws.Range("P10").Value = 0.1283003375911419
$ws.Range("Q10").Value = 0.4159287694924445
$ws.Range("R10").Value = 3.743358925432
$ws.Range("S10").Value = 0.05951724443045241
$ws.Range("T10").Value = 0.05951724443045239

# Row 11
$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Cxcl13"
$ws.Range("C11").Value = "Ackr4"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 3.636334333333334
$ws.Range("H11").Value = 10.909003
$ws.Range("I11").Value = 0.463890006432544
$ws.Range("J11").Value = 0.463890006432544
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.7200953333333334
$ws.Range("N11").Value = 2.160286
$ws.Range("O11").Value = 0.8077233554817153
$ws.Range("P11").Value = 0.8077233554817151
$ws.Range("Q11").Value = 2.618507383873112
$ws.Range("R11").Value = 23.566566454858
$ws.Range("S11").Value = 0.3746947925701289
$ws.Range("T11").Value = 0.3746947925701288

# Row 12
$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Cxcl13"
$ws.Range("C12").Value = "Ackr4"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 3.636334333333334
$ws.Range("H12").Value = 10.909003
$ws.Range("I12").Value = 0.463890006432544
$ws.Range("J12").Value = 0.463890006432544
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.03357866666666667
$ws.Range("N12").Value = 0.100736
$ws.Range("O12").Value = 0.03766483694187069
$ws.Range("P12").Value = 0.03766483694187069
$ws.Range("Q12").Value = 0.1221032584675556
$ws.Range("R12").Value = 1.098929326208
$ws.Range("S12").Value = 0.01747234145124511
$ws.Range("T12").Value = 0.01747234145124511

# Row 13
$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Cxcl13"
$ws.Range("C13").Value = "Ackr4"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 3.636334333333334
$ws.Range("H13").Value = 10.909003
$ws.Range("I13").Value = 0.463890006432544
$ws.Range("J13").Value = 0.463890006432544
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.023457
$ws.Range("N13").Value = 0.070371
$ws.Range("O13").Value = 0.02631146998527222
$ws.Range("P13").Value = 0.02631146998527222
$ws.Range("Q13").Value = 0.08529749445700001
$ws.Range("R13").Value = 0.767677450113
$ws.Range("S13").Value = 0.01220562798071762
$ws.Range("T13").Value = 0.01220562798071761
